$wb = $excel.ActiveWorkbook

# --- Estadisticos 1P ---
$ws1 = $wb.Worksheets.Item("Estadisticos 1P")
$ws1.Range("D4").Value = 0
$ws1.Range("E4").Value = 7
$ws1.Range("H4").Value = 6.5
$ws1.Range("E7").Value = 12
$ws1.Range("F7").Value = 7
$ws1.Range("G7").Value = 36.84
$ws1.Range("H7").Value = 5.5

# --- Estadisticos 2P ---
$ws2 = $wb.Worksheets.Item("Estadisticos 2P")
$ws2.Range("D2").Value = 0
$ws2.Range("E2").Value = 10
$ws2.Range("F2").Value = 26
$ws2.Range("G2").Value = 72.22
$ws2.Range("H2").Value = 5.7
$ws2.Range("D4").Value = 0
$ws2.Range("E4").Value = 5
$ws2.Range("F4").Value = 18
$ws2.Range("G4").Value = 78.26
$ws2.Range("H4").Value = 6.5
$ws2.Range("D5").Value = 0
$ws2.Range("E5").Value = 0
$ws2.Range("F5").Value = 30
$ws2.Range("G5").Value = 100
$ws2.Range("H5").Value = 6.5
$ws2.Range("D6").Value = 0
$ws2.Range("E6").Value = 0
$ws2.Range("F6").Value = 11
$ws2.Range("G6").Value = 100
$ws2.Range("H6").Value = 6.5
$ws2.Range("D7").Value = 0
$ws2.Range("E7").Value = 9
$ws2.Range("F7").Value = 10
$ws2.Range("G7").Value = 52.63
$ws2.Range("H7").Value = 5.5

# --- Estadisticos Final ---
$ws3 = $wb.Worksheets.Item("Estadisticos Final")
$ws3.Range("E2").Value = 10
$ws3.Range("F2").Value = 26
$ws3.Range("G2").Value = 72.22
$ws3.Range("H2").Value = 6.9
$ws3.Range("D4").Value = 0
$ws3.Range("E4").Value = 5
$ws3.Range("F4").Value = 18
$ws3.Range("G4").Value = 78.26
$ws3.Range("H4").Value = 7
$ws3.Range("E5").Value = 0
$ws3.Range("F5").Value = 30
$ws3.Range("G5").Value = 100
$ws3.Range("H5").Value = 7.5
$ws3.Range("E6").Value = 0
$ws3.Range("F6").Value = 11
$ws3.Range("G6").Value = 100
$ws3.Range("H6").Value = 7.2
$ws3.Range("E7").Value = 9
$ws3.Range("F7").Value = 10
$ws3.Range("G7").Value = 52.63
$ws3.Range("H7").Value = 6

# --- Rescatables ---
$ws4 = $wb.Worksheets.Item("Rescatables")
$ws4.Range("A2:G34").ClearContents()
$arr4 = New-Object 'object[,]' 29,7
$arr4[0,0] = 24330051920304
$arr4[0,1] = "ARMAS"
$arr4[0,2] = "SALINAS"
$arr4[0,3] = "JOSE GUSTAVO"
$arr4[0,4] = "Ingles II"
$arr4[0,5] = "2AEV"
$arr4[0,6] = 4
$arr4[1,0] = 24330051920305
$arr4[1,1] = "MORALES"
$arr4[1,2] = "CUAHUA"
$arr4[1,3] = "ANDRES"
$arr4[1,4] = "Ingles II"
$arr4[1,5] = "2AEV"
$arr4[1,6] = 4
$arr4[2,0] = 24330051920113
$arr4[2,1] = "RAMOS"
$arr4[2,2] = "DE LA CRUZ"
$arr4[2,3] = "DEREK"
$arr4[2,4] = "Ingles II"
$arr4[2,5] = "2AEV"
$arr4[2,6] = 4
$arr4[3,0] = 24330051920389
$arr4[3,1] = "RUIZ"
$arr4[3,2] = "MORALES"
$arr4[3,3] = "MAYRIN GUADALUPE"
$arr4[3,4] = "Ingles II"
$arr4[3,5] = "2ALCV"
$arr4[3,6] = 4
$arr4[4,0] = 24330051920247
$arr4[4,1] = "SANCHEZ"
$arr4[4,2] = "PINO"
$arr4[4,3] = "YARETZY NAOMI"
$arr4[4,4] = "Ingles II"
$arr4[4,5] = "2ALCV"
$arr4[4,6] = 4
$arr4[5,0] = 24330051920330
$arr4[5,1] = "VASQUEZ"
$arr4[5,2] = "PEREZ"
$arr4[5,3] = "DANIELA LILI"
$arr4[5,4] = "Ingles II"
$arr4[5,5] = "2ALCV"
$arr4[5,6] = 4
$arr4[6,0] = 24330051920246
$arr4[6,1] = "ZUNO"
$arr4[6,2] = "FLORES"
$arr4[6,3] = "ALIN MARIEL"
$arr4[6,4] = "Ingles II"
$arr4[6,5] = "2ALCV"
$arr4[6,6] = 4
$arr4[7,0] = 23330051920211
$arr4[7,1] = "VAZQUEZ"
$arr4[7,2] = "CARRILLO"
$arr4[7,3] = "DIEGO ARMANDO"
$arr4[7,4] = "Ingles IV"
$arr4[7,5] = "4AEV"
$arr4[7,6] = 4
$arr4[8,0] = 24330051920093
$arr4[8,1] = "ARIAS"
$arr4[8,2] = "SARMIENTO"
$arr4[8,3] = "URIEL ARTURO"
$arr4[8,4] = "Ingles II"
$arr4[8,5] = "2AEV"
$arr4[8,6] = 3
$arr4[9,0] = 24330051920392
$arr4[9,1] = "CERON"
$arr4[9,2] = "GONZALEZ"
$arr4[9,3] = "LEVI SANTIAGO"
$arr4[9,4] = "Ingles II"
$arr4[9,5] = "2AEV"
$arr4[9,6] = 3
$arr4[10,0] = 24330051920098
$arr4[10,1] = "CHICO"
$arr4[10,2] = "BALDERAS"
$arr4[10,3] = "YARETH"
$arr4[10,4] = "Ingles II"
$arr4[10,5] = "2AEV"
$arr4[10,6] = 3
$arr4[11,0] = 24330051920144
$arr4[11,1] = "MUÑOZ"
$arr4[11,2] = "CORONA"
$arr4[11,3] = "JOSE ABEL"
$arr4[11,4] = "Ingles II"
$arr4[11,5] = "2AEV"
$arr4[11,6] = 3
$arr4[12,0] = 24330051920143
$arr4[12,1] = "ROSAS"
$arr4[12,2] = "MEZA"
$arr4[12,3] = "CARLOS ANTONIO"
$arr4[12,4] = "Ingles II"
$arr4[12,5] = "2AEV"
$arr4[12,6] = 3
$arr4[13,0] = 24330051920244
$arr4[13,1] = "BAUTISTA"
$arr4[13,2] = "TORRES"
$arr4[13,3] = "LUZ ELENA"
$arr4[13,4] = "Ingles II"
$arr4[13,5] = "2ALCV"
$arr4[13,6] = 3
$arr4[14,0] = 24330051920238
$arr4[14,1] = "TORRES"
$arr4[14,2] = "PEREZ"
$arr4[14,3] = "ERIKA VALERIA"
$arr4[14,4] = "Ingles II"
$arr4[14,5] = "2ALCV"
$arr4[14,6] = 3
$arr4[15,0] = 22330051920389
$arr4[15,1] = "FLORES"
$arr4[15,2] = "LAGUNA"
$arr4[15,3] = "JOSE ANTONIO"
$arr4[15,4] = "Ingles IV"
$arr4[15,5] = "4AEV"
$arr4[15,6] = 3
$arr4[16,0] = 23330051920332
$arr4[16,1] = "RODRIGUEZ"
$arr4[16,2] = "SUAREZ"
$arr4[16,3] = "SERGIO JOSUE"
$arr4[16,4] = "Ingles IV"
$arr4[16,5] = "4AEV"
$arr4[16,6] = 3
$arr4[17,0] = 23330051920224
$arr4[17,1] = "DORANTES"
$arr4[17,2] = "PORRAS"
$arr4[17,3] = "ROBERTO"
$arr4[17,4] = "Ingles II"
$arr4[17,5] = "2AEV"
$arr4[17,6] = 2
$arr4[18,0] = 24330051920243
$arr4[18,1] = "BERNABE"
$arr4[18,2] = "REYES"
$arr4[18,3] = "JOHANA"
$arr4[18,4] = "Ingles II"
$arr4[18,5] = "2ALCV"
$arr4[18,6] = 2
$arr4[19,0] = 24330051920220
$arr4[19,1] = "GARCIA"
$arr4[19,2] = "CHAPARRO"
$arr4[19,3] = "MAYKA XIMENA"
$arr4[19,4] = "Ingles II"
$arr4[19,5] = "2ALCV"
$arr4[19,6] = 2
$arr4[20,0] = 24330051920226
$arr4[20,1] = "LEYVA"
$arr4[20,2] = "HERNANDEZ"
$arr4[20,3] = "EUNICE GUADALUPE"
$arr4[20,4] = "Ingles II"
$arr4[20,5] = "2ALCV"
$arr4[20,6] = 2
$arr4[21,0] = 24330051920232
$arr4[21,1] = "PELLICO"
$arr4[21,2] = "REYES"
$arr4[21,3] = "BERNARDINA"
$arr4[21,4] = "Ingles II"
$arr4[21,5] = "2ALCV"
$arr4[21,6] = 2
$arr4[22,0] = 23330051920329
$arr4[22,1] = "JIMENEZ"
$arr4[22,2] = "CIRUELO"
$arr4[22,3] = "SABDY"
$arr4[22,4] = "Ingles IV"
$arr4[22,5] = "4AEV"
$arr4[22,6] = 2
$arr4[23,0] = 23330051920324
$arr4[23,1] = "JUAREZ"
$arr4[23,2] = "LIBRADO"
$arr4[23,3] = "ARMANDO GABRIEL"
$arr4[23,4] = "Ingles IV"
$arr4[23,5] = "4AEV"
$arr4[23,6] = 2
$arr4[24,0] = 23330051920203
$arr4[24,1] = "PEREZ"
$arr4[24,2] = "DE JESUS"
$arr4[24,3] = "LUIS FABIAN"
$arr4[24,4] = "Ingles IV"
$arr4[24,5] = "4AEV"
$arr4[24,6] = 2
$arr4[25,0] = 24330051920396
$arr4[25,1] = "MARTINEZ"
$arr4[25,2] = "GONZALEZ"
$arr4[25,3] = "SANTIAGO"
$arr4[25,4] = "Ingles II"
$arr4[25,5] = "2ALCV"
$arr4[25,6] = 1
$arr4[26,0] = 24330051920408
$arr4[26,1] = "MACARIO"
$arr4[26,2] = "SANCHEZ"
$arr4[26,3] = "EMIRETH"
$arr4[26,4] = "Ingles II"
$arr4[26,5] = "2ALCV"
$arr4[26,6] = 1
$arr4[27,0] = 24330051920229
$arr4[27,1] = "MORA"
$arr4[27,2] = "RODRIGUEZ"
$arr4[27,3] = "BRENDA"
$arr4[27,4] = "Ingles II"
$arr4[27,5] = "2ALCV"
$arr4[27,6] = 1
$arr4[28,0] = 24330051920242
$arr4[28,1] = "ZOPIYACTLE"
$arr4[28,2] = "REYES"
$arr4[28,3] = "MARIA DEL CARMEN"
$arr4[28,4] = "Ingles II"
$arr4[28,5] = "2ALCV"
$arr4[28,6] = 1
$ws4.Range("A2:G30").Value = $arr4

